# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F) and "最低票价" (G) figures scraped for
# the 江西-漫展信息 workbook across the relevant sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1905
$ws.Range("G2").Value  = 108
$ws.Range("F4").Value  = 871
$ws.Range("F6").Value  = 54
$ws.Range("F10").Value = 159
$ws.Range("F11").Value = 146
$ws.Range("F13").Value = 4522
$ws.Range("F17").Value = 449
$ws.Range("F18").Value = 15
$ws.Range("F20").Value = 1223
$ws.Range("F21").Value = 2417
$ws.Range("F23").Value = 68
$ws.Range("F24").Value = 47
$ws.Range("F25").Value = 57
$ws.Range("F26").Value = 2238
$ws.Range("F29").Value = 23

# ---- Sheet "演出" ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 37

# ---- Sheet "全部类型" ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 1905
$ws.Range("G2").Value  = 108
$ws.Range("F4").Value  = 871
$ws.Range("F6").Value  = 54
$ws.Range("F10").Value = 159
$ws.Range("F11").Value = 146
$ws.Range("F13").Value = 37
$ws.Range("F14").Value = 4522
$ws.Range("F18").Value = 449
$ws.Range("F19").Value = 15
$ws.Range("F21").Value = 1223
$ws.Range("F22").Value = 2417
$ws.Range("F24").Value = 68
$ws.Range("F25").Value = 47
$ws.Range("F26").Value = 57
$ws.Range("F27").Value = 2238
$ws.Range("F30").Value = 23

$wb.Save()
